$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = "V8F7VA"
$ws.Range("B40").Value = "Cuchilla de limpieza Kyocera"
$ws.Range("C40").Value = "KM 1500 1815 1820, FS 1000 1010 1018 1020 1050"
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 100000
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = 0
$ws.Range("H40").Formula = "=(E40-D40)*G40"
$ws.Range("I40").Formula = "=D40*F40"
$ws.Range("J40").Value = 0
